# Weekly update: insert a new daily price record for
# "Hortaliza, Feria Lagunitas de Puerto Montt - Pepino dulce".
# A new row is inserted at row 33, pushing all subsequent rows down by
# one, and the new row is populated with the latest observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 33 (shifts rows 33:76 down to 34:77).
$ws.Rows("33:33").Insert()

$ws.Range("A33").Value = 4
$ws.Range("B33").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C33").Value = "Los Lagos"
$ws.Range("D33").Value = 44994
$ws.Range("E33").Value = 10
$ws.Range("F33").Value = 100112043
$ws.Range("G33").Value = "Pepino dulce"
$ws.Range("H33").Value = "Cultivar IV Región"
$ws.Range("I33").Value = "Primera"
$ws.Range("J33").Value = 25
$ws.Range("K33").Value = 20000
$ws.Range("L33").Value = 20000
$ws.Range("M33").Value = 20000
$ws.Range("N33").Value = "`$/bandeja 18 kilos"
$ws.Range("O33").Value = "Provincia de Limarí"
$ws.Range("P33").Value = 1111
$ws.Range("Q33").Value = 18
$ws.Range("R33").Value = "Hortaliza"
